$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 100
$ws.Range("I6").Value = 100
$ws.Range("K6").Value = 300
$ws.Range("M6").Value = -188
$ws.Range("H19").Value = 1052.7
$ws.Range("I19").Value = 1150
$ws.Range("K19").Value = 1150
$ws.Range("M19").Value = -975
$ws.Range("H53").Value = 491.18182
$ws.Range("I53").Value = 537.3
$ws.Range("K53").Value = 537.3
$ws.Range("M53").Value = 99.70000000000005
$ws.Range("H74").Value = 2800
$ws.Range("I74").Value = 2800
$ws.Range("K74").Value = 2800
$ws.Range("M74").Value = -1864
$ws.Range("H77").Value = 2800
$ws.Range("I77").Value = 2800
$ws.Range("K77").Value = 14000
$ws.Range("M77").Value = -9320
$ws.Range("H80").Value = 753.2727
$ws.Range("I80").Value = 1088
$ws.Range("J80").Value = 562
$ws.Range("K80").Value = 3264
$ws.Range("L80").Value = 1686
$ws.Range("M80").Value = -2266
$ws.Range("N80").Value = -3682
$ws.Range("H83").Value = 753.2727
$ws.Range("I83").Value = 1088
$ws.Range("J83").Value = 562
$ws.Range("K83").Value = 9792
$ws.Range("L83").Value = 5058
$ws.Range("M83").Value = -4800
$ws.Range("N83").Value = -15042
$ws.Range("H137").Value = 2234.4614
$ws.Range("I137").Value = 2370.6667
$ws.Range("K137").Value = 7112.000100000001
$ws.Range("M137").Value = -4562.000100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 35989.6
$ws.Range("J62").Value = 35989.6
$ws.Range("L62").Value = 35989.6
$ws.Range("N62").Value = -37237.6
$ws.Range("H65").Value = 35989.6
$ws.Range("J65").Value = 35989.6
$ws.Range("L65").Value = 107968.8
$ws.Range("N65").Value = -114208.8
$ws.Range("H122").Value = 1978
$ws.Range("I122").Value = 1975
$ws.Range("K122").Value = 5925
$ws.Range("M122").Value = -3475
$ws.Range("H132").Value = 6644
$ws.Range("I132").Value = 6644
$ws.Range("K132").Value = 19932
$ws.Range("M132").Value = -17402

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 3800
$ws.Range("J9").Value = 3800
$ws.Range("L9").Value = 3800
$ws.Range("N9").Value = -4136
$ws.Range("H22").Value = 584.8889
$ws.Range("I22").Value = 537.8570999999999
$ws.Range("J22").Value = 749.5
$ws.Range("K22").Value = 537.8570999999999
$ws.Range("L22").Value = 749.5
$ws.Range("M22").Value = -364.8570999999999
$ws.Range("N22").Value = -1095.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 4750
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H3").Value = 5000000
$ws.Range("I3").Value = 5000000
$ws.Range("K3").Value = 5000000
$ws.Range("M3").Value = -4999887
$ws.Range("H4").Value = 366.66666
$ws.Range("I4").Value = 366.66666
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 366.66666
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -254.66666
$ws.Range("N4").ClearContents()
$ws.Range("H14").Value = 2666.6667
$ws.Range("I14").Value = 2300
$ws.Range("J14").Value = 3400
$ws.Range("K14").Value = 2300
$ws.Range("L14").Value = 3400
$ws.Range("M14").Value = -2130
$ws.Range("N14").Value = -3740
$ws.Range("H15").Value = 4454.5
$ws.Range("J15").Value = 4454.5
$ws.Range("L15").Value = 4454.5
$ws.Range("M15").Value = -4794.5
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1541
$ws.Range("J4").Value = 1649.2
$ws.Range("L4").Value = 4947.6
$ws.Range("N4").Value = -5171.6
$ws.Range("H41").Value = 1275
$ws.Range("J41").Value = 1275
$ws.Range("L41").Value = 3825
$ws.Range("N41").Value = -4501
$ws.Range("H64").Value = 222
$ws.Range("I64").Value = 222
$ws.Range("K64").Value = 666
$ws.Range("M64").Value = -396
$ws.Range("H67").Value = 222
$ws.Range("I67").Value = 222
$ws.Range("K67").Value = 666
$ws.Range("M67").Value = 270
$ws.Range("H92").Value = 2000
$ws.Range("J92").Value = 2000
$ws.Range("L92").Value = 6000
$ws.Range("N92").Value = -8496
$ws.Range("H131").Value = 4020.9092
$ws.Range("I131").Value = 3410
$ws.Range("J131").Value = 4250
$ws.Range("K131").Value = 10230
$ws.Range("L131").Value = 12750
$ws.Range("M131").Value = -5190
$ws.Range("N131").Value = -22830

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 45000000
$ws.Range("I35").Value = 45000000
$ws.Range("K35").Value = 45000000
$ws.Range("M35").Value = -44999702
$ws.Range("H63").Value = 50000
$ws.Range("J63").Value = 50000
$ws.Range("L63").Value = 50000
$ws.Range("N63").Value = -51372
$ws.Range("H66").Value = 50000
$ws.Range("J66").Value = 50000
$ws.Range("L66").Value = 150000
$ws.Range("N66").Value = -156864
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H122").Value = 2000
$ws.Range("J122").Value = 2000
$ws.Range("L122").Value = 6000
$ws.Range("N122").Value = -10900
$ws.Range("H132").Value = 2508
$ws.Range("I132").Value = 2635
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 7905
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -5375
$ws.Range("N132").Value = -11060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2642.3333
$ws.Range("I16").Value = 2722.625
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 2722.625
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -2552.625
$ws.Range("N16").Value = -2340
$ws.Range("H40").Value = 2998.5
$ws.Range("I40").Value = 2998.5
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2998.5
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2862.5
$ws.Range("N40").ClearContents()
$ws.Range("H43").Value = 24335
$ws.Range("I43").Value = 23670
$ws.Range("K43").Value = 23670
$ws.Range("M43").Value = -23477
$ws.Range("H45").Value = 20000000
$ws.Range("I45").Value = 20000000
$ws.Range("K45").Value = 20000000
$ws.Range("M45").Value = -19999593
$ws.Range("H132").Value = 52502
$ws.Range("I132").Value = 5004
$ws.Range("K132").Value = 15012
$ws.Range("M132").Value = -12482

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 493
$ws.Range("I81").Value = 493
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 986
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = 75
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 493
$ws.Range("I84").Value = 493
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 4930
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 374
$ws.Range("N84").ClearContents()
$ws.Range("H132").Value = 899.6667
$ws.Range("I132").Value = 899.5
$ws.Range("K132").Value = 2698.5
$ws.Range("M132").Value = -168.5
$ws.Range("H136").Value = 128528.125
$ws.Range("I136").Value = 3075.2
$ws.Range("K136").Value = 9225.599999999999
$ws.Range("M136").Value = -6675.599999999999

Write-Host "Applied all Rafflesia_Profits cell updates"